# Applies the 2025-04-08 vault-backup edit to fangfeiyue.docx
#
# 1) Merge the two runs around the stray "_GoBack" bookmark in the
#    "分销团队..." sentence back into a single run (bookmark removed).
# 2) "合理" -> "良好" in the after-sales paragraph.
# 3) Rewrite the third-party logistics margin sentence.
# 4) Split the "有赞寄件" sentence into two runs, re-inserting the
#    "_GoBack" bookmark at the split point.

$d = $word.ActiveDocument

# --- Change 1 ------------------------------------------------------
# The two runs already read as one continuous text stream across the
# bookmark, so a Find/Replace with identical find/replace text merges
# them into a single run and drops the now-redundant bookmark.
$d.Content.Find.Execute(
    "分销团队，虚线管理云分销团队、并负责TOP级插件分销员的架构优化及稳定性治理",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "分销团队，虚线管理云分销团队、并负责TOP级插件分销员的架构优化及稳定性治理",
    2) | Out-Null

# --- Change 2 --------------------------------------------------------
$d.Content.Find.Execute(
    "为消费者提供合理的售后履约服务，并通过动态定价，来赚取商家服务费。",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "为消费者提供良好的售后履约服务，并通过动态定价，来赚取商家服务费。",
    2) | Out-Null

# --- Change 3 --------------------------------------------------------
$d.Content.Find.Execute(
    "通过物流商运营、定价运营，赚取与三方物流间的差价。",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "为消费者提供具有性价比的三方物流服务，并通过物流商运营、定价运营等，赚取与三方间的差价。",
    2) | Out-Null

# --- Change 4 --------------------------------------------------------
# Find the first half of the "有赞寄件" sentence, collapse the found
# range to its end, and drop a "_GoBack" bookmark there -- this splits
# the original single run into two runs with the bookmark in between,
# matching the diff.
$r = $d.Content
$r.Find.Execute(
    "有赞寄件：通过承接正向的交易发货能力，",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "", 0) | Out-Null
$r.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r) | Out-Null
